$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.307.22"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$ws.Range("D3").Value = "3.590.03"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.32"
$ws.Range("E5").Value = "  -2.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.49"
$ws.Range("E6").Value = "  -3.30%  "

# Row 7
$ws.Range("D7").Value = "3.585.48"
$ws.Range("E7").Value = "  -1.23%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.623"
$ws.Range("E8").Value = "  -3.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").Value = "  +0.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.655"
$ws.Range("E11").Value = "  -3.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.41"

# Row 13
$ws.Range("E13").Value = "  -2.40%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.59"
$ws.Range("E14").Value = "  -3.51%  "

# Row 15
$ws.Range("D15").Value = "4.161.07"
$ws.Range("E15").Value = "  -1.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.72"

# Row 17
$ws.Range("D17").Value = "3.591.08"
$ws.Range("E17").Value = "  -1.24%  "

# Row 18
$ws.Range("D18").Value = "70.204.90"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.58"
$ws.Range("E19").Value = "  -1.32%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.120"
$ws.Range("E20").Value = "  -1.22%  "

# Row 21
$ws.Range("E21").Value = "  -2.47%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.32"
$ws.Range("E22").Value = "  +1.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.66"
$ws.Range("E23").Value = "  +1.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.99"
$ws.Range("E24").Value = "  -4.25%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.42"
$ws.Range("E25").Value = "  +6.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.41"
$ws.Range("E26").Value = "  -1.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.99"
$ws.Range("E28").Value = "  -5.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("E29").Value = "  -1.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("E30").Value = "  -2.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.79"
$ws.Range("E31").Value = "  -3.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.84"
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("E34").Value = "  -5.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "573.93"
$ws.Range("E35").Value = "  -6.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.18"
$ws.Range("E36").Value = "  +10.92%  "

# Row 37
$ws.Range("E37").Value = "  -3.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.406"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("E39").Value = "  +0.02%  "

# Row 40
$ws.Range("E40").Value = "  -4.97%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  -2.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.20"
$ws.Range("E42").Value = "  +1.75%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -8.39%  "

# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  +9.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.10"
$ws.Range("E45").Value = "  -2.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0455"
$ws.Range("E46").Value = "  -0.69%  "

# Row 47
$ws.Range("D47").Value = "3.207.84"
$ws.Range("E47").Value = "  -3.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.48"
$ws.Range("E48").Value = "  -1.86%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.137"
$ws.Range("E49").Value = "  -1.90%  "

# Row 50
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.03%  "

# Row 51
$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.47"
$ws.Range("E51").Value = "  +22.02%  "
